# Regenerate save_data to use K (strikeouts) instead of Strike# (pitches-that-were-strikes count).
# This recalculates/rewrites the "K" column (column G) values for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 2
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 2
    12 = 2
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 2
    23 = 1
    24 = 2
    25 = 1
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 2
    31 = 0
    32 = 1
    33 = 0
    34 = 0
    35 = 0
    36 = 2
    37 = 2
    38 = 0
    39 = 1
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 2
    45 = 1
    46 = 2
    47 = 1
    48 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
